$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.220.81'
$ws.Range('E2').Value = '  -1.33%  '

$ws.Range('D3').Value = '3.534.91'
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.26'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.76%  '

$ws.Range('D7').Value = '3.533.55'

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.480'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.137'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '8.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.97%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.413'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.47%  '

$ws.Range('D13').Value = '4.132.13'
$ws.Range('E13').Value = '  +0.38%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000208'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.48%  '

$ws.Range('E15').Value = '  -4.74%  '

$ws.Range('D16').Value = '3.527.81'
$ws.Range('E16').Value = '  +0.11%  '

$ws.Range('D17').Value = '66.274.16'
$ws.Range('E17').Value = '  -1.23%  '

$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.99%  '

$ws.Range('E20').Value = '  -3.09%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '425.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.80%  '

$ws.Range('E23').Value = '  -1.43%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.68%  '

$ws.Range('D25').Value = '3.671.38'
$ws.Range('E25').Value = '  +0.21%  '

$ws.Range('E26').Value = '  -0.01%  '

$ws.Range('E27').Value = '  -1.80%  '

$ws.Range('E28').Value = '  -5.03%  '

$ws.Range('E29').Value = '  -2.85%  '

$ws.Range('E30').Value = '  -1.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.13%  '

$ws.Range('E32').Value = '  -3.09%  '

$ws.Range('E33').Value = '  -6.26%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.29%  '

$ws.Range('D35').Value = '3.523.70'
$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('E36').Value = '  -0.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.75'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.10%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.36%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.61'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.20%  '

$ws.Range('E40').Value = '  -0.13%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '171.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.64%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0860'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.89%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.894'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.97%  '

$ws.Range('E45').Value = '  -9.22%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.24'
$ws.Range('D46').Style = 'Normal'

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.46%  '

$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.14%  '

$ws.Range('E49').Value = '  -1.12%  '

$ws.Range('E50').Value = '  -4.10%  '

$ws.Range('E51').Value = '  -4.15%  '
